$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5298.1665
$ws.Range("J40").Value = 1461.8
$ws.Range("L40").Value = 1461.8
$ws.Range("N40").Value = -1811.8
$ws.Range("H86").Value = 3085.2222
$ws.Range("I86").Value = 2176.3635
$ws.Range("K86").Value = 2176.3635
$ws.Range("M86").Value = -1053.3635
$ws.Range("H88").Value = 2626.25
$ws.Range("I88").Value = 2252.5
$ws.Range("J88").Value = 3000
$ws.Range("K88").Value = 2252.5
$ws.Range("L88").Value = 3000
$ws.Range("M88").Value = -1846.5
$ws.Range("N88").Value = -3812
$ws.Range("H89").Value = 3085.2222
$ws.Range("I89").Value = 2176.3635
$ws.Range("K89").Value = 10881.8175
$ws.Range("M89").Value = -5265.817499999999
$ws.Range("H91").Value = 2626.25
$ws.Range("I91").Value = 2252.5
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 2252.5
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = -848.5
$ws.Range("N91").Value = -5808
$ws.Range("H138").Value = 6669714
$ws.Range("J138").Value = 7815792.5
$ws.Range("L138").Value = 23447377.5
$ws.Range("N138").Value = -23457657.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 538.3043
$ws.Range("I2").Value = 308.6111
$ws.Range("K2").Value = 308.6111
$ws.Range("M2").Value = -195.6111
$ws.Range("H32").Value = 7482.8315
$ws.Range("I32").Value = 3204.782
$ws.Range("J32").Value = 27111.53
$ws.Range("K32").Value = 3204.782
$ws.Range("L32").Value = 27111.53
$ws.Range("M32").Value = -2917.782
$ws.Range("N32").Value = -27685.53
$ws.Range("H61").Value = 5203.5884
$ws.Range("I61").Value = 4265.636
$ws.Range("K61").Value = 4265.636
$ws.Range("M61").Value = -4053.636
$ws.Range("H74").Value = 72057.375
$ws.Range("I74").Value = 112692.89
$ws.Range("J74").Value = 19811.715
$ws.Range("K74").Value = 112692.89
$ws.Range("L74").Value = 19811.715
$ws.Range("M74").Value = -111818.89
$ws.Range("N74").Value = -21559.715
$ws.Range("H77").Value = 72057.375
$ws.Range("I77").Value = 112692.89
$ws.Range("J77").Value = 19811.715
$ws.Range("K77").Value = 563464.45
$ws.Range("L77").Value = 99058.575
$ws.Range("M77").Value = -559096.45
$ws.Range("N77").Value = -107794.575
$ws.Range("H97").Value = 2002.6923
$ws.Range("I97").Value = 1523.75
$ws.Range("J97").Value = 7750
$ws.Range("K97").Value = 1523.75
$ws.Range("L97").Value = 7750
$ws.Range("M97").Value = -1027.75
$ws.Range("N97").Value = -8742
$ws.Range("H102").Value = 4862
$ws.Range("I102").Value = 5077.5
$ws.Range("K102").Value = 5077.5
$ws.Range("M102").Value = -3455.5
$ws.Range("H116").Value = 538.3043
$ws.Range("I116").Value = 308.6111
$ws.Range("K116").Value = 308.6111
$ws.Range("M116").Value = 1985.3889
$ws.Range("H132").Value = 3342.8628
$ws.Range("I132").Value = 3222.2
$ws.Range("J132").Value = 3781.6365
$ws.Range("K132").Value = 9666.599999999999
$ws.Range("L132").Value = 11344.9095
$ws.Range("M132").Value = -7136.599999999999
$ws.Range("N132").Value = -16404.9095
$ws.Range("H136").Value = 5203.5884
$ws.Range("I136").Value = 4265.636
$ws.Range("K136").Value = 12796.908
$ws.Range("M136").Value = -10246.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 538.3043
$ws.Range("I3").Value = 308.6111
$ws.Range("K3").Value = 308.6111
$ws.Range("M3").Value = -194.6111
$ws.Range("H10").Value = 1000
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1000
$ws.Range("L10").ClearContents()
$ws.Range("N10").Value = 0
$ws.Range("M10").Value = -860
$ws.Range("H99").Value = 148313.42
$ws.Range("I99").Value = 113899.78
$ws.Range("J99").Value = 210258
$ws.Range("K99").Value = 113899.78
$ws.Range("L99").Value = 210258
$ws.Range("M99").Value = -112401.78
$ws.Range("N99").Value = -213254
$ws.Range("H105").Value = 2256.2368
$ws.Range("I105").Value = 2188.6453
$ws.Range("K105").Value = 2188.6453
$ws.Range("M105").Value = -441.6453000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35305.582
$ws.Range("I31").Value = 48164.684
$ws.Range("J31").Value = 3872.2222
$ws.Range("K31").Value = 48164.684
$ws.Range("L31").Value = 3872.2222
$ws.Range("M31").Value = -47869.684
$ws.Range("N31").Value = -4462.2222
$ws.Range("H34").Value = 35305.582
$ws.Range("I34").Value = 48164.684
$ws.Range("J34").Value = 3872.2222
$ws.Range("K34").Value = 48164.684
$ws.Range("L34").Value = 3872.2222
$ws.Range("M34").Value = -47962.684
$ws.Range("N34").Value = -4276.2222
$ws.Range("H86").Value = 5459.387
$ws.Range("I86").Value = 5309.407
$ws.Range("K86").Value = 5309.407
$ws.Range("M86").Value = -4186.407
$ws.Range("H89").Value = 5459.387
$ws.Range("I89").Value = 5309.407
$ws.Range("K89").Value = 26547.035
$ws.Range("M89").Value = -20931.035
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0
$ws.Range("H133").Value = 59999
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H140").Value = 115000
$ws.Range("J140").Value = 115000
$ws.Range("L140").Value = 115000
$ws.Range("N140").Value = -125360
$ws.Range("H141").Value = 733442
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 733442
$ws.Range("K141").Value = 0
$ws.Range("L141").ClearContents()
$ws.Range("M141").Value = 733442
$ws.Range("N141").Value = -743802

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 18531518
$ws.Range("I18").Value = 22225822
$ws.Range("K18").Value = 22225822
$ws.Range("M18").Value = -22225529
$ws.Range("H113").Value = 3719.9092
$ws.Range("I113").Value = 3691.9
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3691.9
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -1521.9
$ws.Range("N113").Value = -8340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 10257.462
$ws.Range("J20").Value = 10334.7
$ws.Range("L20").Value = 10334.7
$ws.Range("N20").Value = -10786.7
$ws.Range("H55").Value = 229.15384
$ws.Range("J55").Value = 123.333336
$ws.Range("L55").Value = 123.333336
$ws.Range("N55").Value = -469.333336
$ws.Range("H74").Value = 20000
$ws.Range("I74").Value = 20000
$ws.Range("K74").Value = 20000
$ws.Range("M74").Value = -19002
$ws.Range("H77").Value = 20000
$ws.Range("I77").Value = 20000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55008
$ws.Range("H93").Value = 1886.5333
$ws.Range("I93").Value = 1930.3077
$ws.Range("J93").Value = 1602
$ws.Range("K93").Value = 1930.3077
$ws.Range("L93").Value = 1602
$ws.Range("M93").Value = -682.3077000000001
$ws.Range("N93").Value = -4098
$ws.Range("H100").Value = 1976.1111
$ws.Range("I100").Value = 1359
$ws.Range("J100").Value = 2747.5
$ws.Range("K100").Value = 1359
$ws.Range("L100").Value = 2747.5
$ws.Range("M100").Value = -818
$ws.Range("N100").Value = -3829.5
$ws.Range("H122").Value = 479993.06
$ws.Range("I122").Value = 669997
$ws.Range("K122").Value = 2009991
$ws.Range("M122").Value = -2007541
$ws.Range("H136").Value = 3922.125
$ws.Range("I136").Value = 3719.25
$ws.Range("K136").Value = 11157.75
$ws.Range("M136").Value = -8607.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 699.9524
$ws.Range("I113").Value = 756.4
$ws.Range("K113").Value = 2269.2
$ws.Range("M113").Value = -99.19999999999982
$ws.Range("H126").Value = 12570.4
$ws.Range("I126").Value = 12570.4
$ws.Range("K126").Value = 37711.2
$ws.Range("M126").Value = -35241.2

Write-Host "Applied all updates"